# Weekly fruit/vegetable data refresh: a new daily price record (2022-08-09,
# serial 44782) is inserted as row 329, which pushes every existing record
# from row 329 down through row 405 one row further down (330-406). The
# record that used to be at the very end (row 405) becomes the new last row
# (406).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 329; Excel shifts rows 329:405 down to 330:406.
$ws.Rows("329:329").Insert()

# Populate the newly inserted row 329 with the new record's data.
$ws.Range("A329").Value = 9
$ws.Range("B329").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C329").Value = "Metropolitana"
$ws.Range("D329").Value = 44782
$ws.Range("E329").Value = 13
$ws.Range("F329").Value = 100112032
$ws.Range("G329").Value = "Zapallo italiano"
$ws.Range("H329").Value = "Sin especificar"
$ws.Range("I329").Value = "Primera"
$ws.Range("J329").Value = 70
$ws.Range("K329").Value = 21000
$ws.Range("L329").Value = 23000
$ws.Range("M329").Value = 22000
$ws.Range("N329").Value = "$/caja 50 unidades"
$ws.Range("O329").Value = "Región de Arica y Parinacota"
$ws.Range("P329").Value = 440
$ws.Range("Q329").Value = 50
$ws.Range("R329").Value = "Hortaliza"
